$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''37.182.17'
$ws.Range("E2").Value = '  +1.59%  '

# Row 3
$ws.Range("D3").Value = '''2.010.01'
$ws.Range("E3").Value = '  +2.40%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '''246.94'
$ws.Range("E5").Value = '  +1.02%  '

# Row 6
$ws.Range("E6").Value = '  +2.51%  '

# Row 7
$ws.Range("D7").Value = '''59.74'
$ws.Range("E7").Value = '  -1.48%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("E9").Value = '  +3.33%  '

# Row 10
$ws.Range("D10").Value = '''0.0812'
$ws.Range("E10").Value = '  +1.80%  '

# Row 11
$ws.Range("E11").Value = '  +0.79%  '

# Row 12
$ws.Range("D12").Value = '''15.09'
$ws.Range("E12").Value = '  +5.44%  '

# Row 13
$ws.Range("D13").Value = '''22.33'
$ws.Range("E13").Value = '  +2.30%  '

# Row 14
$ws.Range("D14").Value = '''2.302.59'
$ws.Range("E14").Value = '  +2.27%  '

# Row 15
$ws.Range("D15").Value = '''0.846'
$ws.Range("E15").Value = '  +0.77%  '

# Row 16
$ws.Range("E16").Value = '  +3.14%  '

# Row 17
$ws.Range("D17").Value = '''2.012.51'
$ws.Range("E17").Value = '  +2.51%  '

# Row 18
$ws.Range("D18").Value = '''37.112.52'
$ws.Range("E18").Value = '  +1.49%  '

# Row 19
$ws.Range("D19").Value = '''70.31'
$ws.Range("E19").Value = '  +0.41%  '

# Row 20
$ws.Range("E20").Value = '  +1.13%  '

# Row 21
$ws.Range("D21").Value = '''5.22'
$ws.Range("E21").Value = '  +2.34%  '

# Row 22
$ws.Range("D22").Value = '''230.81'
$ws.Range("E22").Value = '  +0.10%  '

# Row 23
$ws.Range("E23").Value = '  +0.01%  '

# Row 24
$ws.Range("E24").Value = '  +0.51%  '

# Row 25
$ws.Range("E25").Value = '  +0.40%  '

# Row 26
$ws.Range("D26").Value = '''9.46'
$ws.Range("E26").Value = '  +2.49%  '

# Row 27
$ws.Range("E27").Value = '  +2.18%  '

# Row 28
$ws.Range("D28").Value = '''0.140'
$ws.Range("E28").Value = '  -3.40%  '

# Row 29
$ws.Range("D29").Value = '''19.71'
$ws.Range("E29").Value = '  +1.32%  '

# Row 30
$ws.Range("E30").Value = '  +12.34%  '

# Row 31
$ws.Range("E31").Value = '  +1.25%  '

# Row 32
$ws.Range("E32").Value = '  +0.81%  '

# Row 33
$ws.Range("D33").Value = '''0.0656'
$ws.Range("E33").Value = '  +5.99%  '

# Row 34
$ws.Range("E34").Value = '  -0.11%  '

# Row 35
$ws.Range("D35").Value = '''2.46'
$ws.Range("E35").Value = '  +7.97%  '

# Row 36
$ws.Range("D36").Value = '''3.47'
$ws.Range("E36").Value = '  -4.21%  '

# Row 37
$ws.Range("E37").Value = '  +0.05%  '

# Row 38
$ws.Range("E38").Value = '  +2.23%  '

# Row 39
$ws.Range("D39").Value = '''5.33'
$ws.Range("E39").Value = '  -4.48%  '

# Row 40
$ws.Range("D40").Value = '''0.0982'
$ws.Range("E40").Value = '  -0.49%  '

# Row 41
$ws.Range("E41").Value = '  +0.85%  '

# Row 42
$ws.Range("E42").Value = '  +0.77%  '

# Row 43
$ws.Range("E43").Value = '  +1.34%  '

# Row 44
$ws.Range("D44").Value = '''16.63'
$ws.Range("E44").Value = '  +2.25%  '

# Row 45
$ws.Range("D45").Value = '''92.13'
$ws.Range("E45").Value = '  +3.87%  '

# Row 46
$ws.Range("D46").Value = '''1.372.85'
$ws.Range("E46").Value = '  +0.41%  '

# Row 47
$ws.Range("E47").Value = '  +1.33%  '

# Row 48
$ws.Range("D48").Value = '''7.41'
$ws.Range("E48").Value = '  +3.29%  '

# Row 49
$ws.Range("D49").Value = '''2.08'
$ws.Range("E49").Value = '  +12.75%  '

# Row 50
$ws.Range("B50").Value = '''MXToken'
$ws.Range("C50").Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '''2.84'
$ws.Range("E50").Value = '  -0.02%  '

# Row 51
$ws.Range("B51").Value = '''MultiversX'
$ws.Range("C51").Value = '''https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '''46.68'
$ws.Range("E51").Value = '  +5.25%  '
